$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "{'correlação_pearson': np.float64(-0.1), 'p_valor': np.float64(0.12)}"
$ws.Range("B3").Value = "{'qui2': np.float64(0.0), 'p_valor': np.float64(1.0), 'graus_liberdade': 1, 'tabela': {1.0: {1.0: 48, 2.0: 54}, 2.0: {1.0: 55, 2.0: 60}}}"
$ws.Range("B4").Value = "{'qui2': np.float64(1.39), 'p_valor': np.float64(0.24), 'graus_liberdade': 1, 'tabela': {1.0: {1.0: 60, 2.0: 54}, 2.0: {1.0: 45, 2.0: 58}}}"
$ws.Range("B5").Value = "{1.0: 6.21, 2.0: 5.2, 3.0: 5.62}"
$ws.Range("B6").Value = "{0: {'Q5': 1.54, 'Q6': 1.38, 'Q7': 4.32, 'Q8': 2.65, 'Q9': 13.19, 'Q10': 1.62, 'Q11': 2.82, 'Q12': 2.13, 'Q13': 4.81, 'Q14': 1.76, 'Q15': 1.34, 'Q16': 1.46, 'Q17': 1.57, 'Q18': 1.51, 'Q19': 1.76, 'Q20': 1.38, 'Q21': 1.56, 'Q22': 1.69, 'Q23': 1.68, 'Q24': 1.43, 'Q25': 1.6, 'Q26': 1.65}, 1: {'Q5': 1.46, 'Q6': 1.45, 'Q7': 4.46, 'Q8': 1.66, 'Q9': 12.71, 'Q10': 1.32, 'Q11': 4.21, 'Q12': 1.64, 'Q13': 5.99, 'Q14': 1.44, 'Q15': 1.28, 'Q16': 1.21, 'Q17': 1.48, 'Q18': 1.41, 'Q19': 1.32, 'Q20': 1.66, 'Q21': 1.38, 'Q22': 1.48, 'Q23': 1.53, 'Q24': 1.59, 'Q25': 1.4, 'Q26': 2.16}, 2: {'Q5': 1.56, 'Q6': 1.78, 'Q7': 4.16, 'Q8': 1.92, 'Q9': 13.97, 'Q10': 1.8, 'Q11': 5.19, 'Q12': 2.23, 'Q13': 6.09, 'Q14': 1.41, 'Q15': 1.73, 'Q16': 1.84, 'Q17': 1.39, 'Q18': 1.45, 'Q19': 1.58, 'Q20': 1.67, 'Q21': 1.45, 'Q22': 1.42, 'Q23': 1.38, 'Q24': 1.52, 'Q25': 1.61, 'Q26': 2.2}}"
